# Apply the updated crypto price / volume(1h) figures produced by the
# scheduled GitHub Actions data refresh.
#
# The "Price" (D) and "Volume(1h)" (E) columns are stored as plain text in
# this workbook (they were written as text, not numbers, by the scraper).
# Several "Price" values look like ordinary decimals (e.g. 1.003, 317.27), so
# assigning them straight to .Value would make Excel auto-convert them to
# numbers. A leading apostrophe forces Excel to keep them as text, matching
# the original data type. "Volume(1h)" values (e.g. "  +1.93%  ") and
# thousands-grouped "Price" values (e.g. 27.962.39) never parse as numbers,
# so they need no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.962.39'; E = '  +1.93%  ' },
    @{ Row = 3; D = '1.904.45'; E = '  +2.34%  ' },
    @{ Row = 4; D = '''1.003'; E = '  -0.80%  ' },
    @{ Row = 5; D = '''317.27'; E = '  +1.93%  ' },
    @{ Row = 6; D = $null; E = '  -0.69%  ' },
    @{ Row = 7; D = '''0.4815'; E = '  +0.88%  ' },
    @{ Row = 8; D = $null; E = '  -0.16%  ' },
    @{ Row = 9; D = '''0.07355'; E = '  +0.50%  ' },
    @{ Row = 10; D = '''0.9307'; E = '  -0.06%  ' },
    @{ Row = 11; D = '''20.76'; E = '  -0.13%  ' },
    @{ Row = 12; D = '''0.07748'; E = '  -0.57%  ' },
    @{ Row = 13; D = '1.883.81'; E = '  +1.19%  ' },
    @{ Row = 14; D = '''5.479'; E = $null },
    @{ Row = 15; D = '''6.634'; E = '  +1.26%  ' },
    @{ Row = 16; D = '''91.57'; E = '  +1.62%  ' },
    @{ Row = 17; D = '''1.005'; E = '  -0.71%  ' },
    @{ Row = 18; D = '''0.000008868'; E = '  +0.52%  ' },
    @{ Row = 19; D = $null; E = '  -0.67%  ' },
    @{ Row = 20; D = '27.996.94'; E = '  +2.01%  ' },
    @{ Row = 21; D = $null; E = '  +0.37%  ' },
    @{ Row = 22; D = '''5.134'; E = '  +0.73%  ' },
    @{ Row = 23; D = '2.144.99'; E = '  +2.09%  ' },
    @{ Row = 24; D = $null; E = '  +1.83%  ' },
    @{ Row = 25; D = '''155.95'; E = '  +0.45%  ' },
    @{ Row = 26; D = '''1.910'; E = '  -1.55%  ' },
    @{ Row = 27; D = '''18.48'; E = '  +0.09%  ' },
    @{ Row = 28; D = '''2.109'; E = '  +4.95%  ' },
    @{ Row = 29; D = '''117.21'; E = '  +1.50%  ' },
    @{ Row = 30; D = '''4.962'; E = '  +0.22%  ' },
    @{ Row = 31; D = '''0.08945'; E = '  +0.54%  ' },
    @{ Row = 32; D = '''3.250'; E = '  -2.46%  ' },
    @{ Row = 33; D = '''1.248'; E = '  +3.55%  ' },
    @{ Row = 34; D = '''0.7689'; E = '  +2.07%  ' },
    @{ Row = 35; D = '''4.658'; E = '  +1.36%  ' },
    @{ Row = 36; D = '''2.587'; E = '  -4.60%  ' },
    @{ Row = 37; D = $null; E = '  +0.25%  ' },
    @{ Row = 38; D = $null; E = '  -1.77%  ' },
    @{ Row = 39; D = '''0.5497'; E = '  -1.09%  ' },
    @{ Row = 40; D = '''0.05286'; E = '  +0.26%  ' },
    @{ Row = 41; D = '''2.993'; E = '  +0.12%  ' },
    @{ Row = 42; D = $null; E = '  -0.86%  ' },
    @{ Row = 43; D = '''0.1527'; E = '  +0.33%  ' },
    @{ Row = 44; D = '''8.467'; E = '  -1.88%  ' },
    @{ Row = 45; D = '''110.16'; E = '  +6.91%  ' },
    @{ Row = 46; D = $null; E = '  +0.70%  ' },
    @{ Row = 47; D = '''0.4817'; E = '  -1.25%  ' },
    @{ Row = 48; D = '''1.003'; E = '  -0.71%  ' },
    @{ Row = 49; D = '''1.640'; E = '  -1.70%  ' },
    @{ Row = 50; D = '''67.70'; E = '  +0.42%  ' },
    @{ Row = 51; D = '''0.06076'; E = '  -0.26%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
